{"js": "// Ordered list of (old, new) text replacements, in document order, taken\n// from the OOXML diff: the date heading plus the 25 two-digit / one-digit\n// division problems in the table. One value (\"95\u00f75=\") occurs twice in the\n// source with two different replacements, so replacements are applied by\n// walking each unique search hit list in document order.\nconst replacements = [\n  [\"2026-01-06 Tuesday\", \"2026-01-07 Wednesday\"],\n  [\"84\u00f75=\", \"30\u00f75=\"],\n  [\"95\u00f75=\", \"98\u00f72=\"],\n  [\"53\u00f76=\", \"57\u00f78=\"],\n  [\"78\u00f79=\", \"31\u00f78=\"],\n  [\"13\u00f74=\", \"91\u00f72=\"],\n  [\"64\u00f79=\", \"42\u00f75=\"],\n  [\"57\u00f79=\", \"64\u00f77=\"],\n  [\"97\u00f73=\", \"59\u00f75=\"],\n  [\"95\u00f75=\", \"57\u00f74=\"],\n  [\"19\u00f78=\", \"59\u00f75=\"],\n  [\"36\u00f74=\", \"81\u00f72=\"],\n  [\"34\u00f77=\", \"92\u00f75=\"],\n  [\"61\u00f79=\", \"86\u00f77=\"],\n  [\"20\u00f76=\", \"78\u00f73=\"],\n  [\"77\u00f73=\", \"18\u00f74=\"],\n  [\"72\u00f76=\", \"88\u00f74=\"],\n  [\"41\u00f75=\", \"99\u00f74=\"],\n  [\"68\u00f78=\", \"17\u00f77=\"],\n  [\"85\u00f79=\", \"69\u00f79=\"],\n  [\"48\u00f74=\", \"35\u00f72=\"],\n  [\"39\u00f74=\", \"35\u00f74=\"],\n  [\"38\u00f76=\", \"54\u00f78=\"],\n  [\"16\u00f78=\", \"60\u00f78=\"],\n  [\"25\u00f79=\", \"98\u00f73=\"],\n  [\"25\u00f72=\", \"93\u00f76=\"],\n];\n\nconst body = context.document.body;\n\n// Search once per distinct \"old\" string, loading every match (matches come\n// back in document order), then consume hits for repeated strings in order.\nconst uniqueOld = [...new Set(replacements.map((p) => p[0]))];\nconst searchResults = {};\nfor (const oldText of uniqueOld) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searchResults[oldText] = found;\n}\nawait context.sync();\n\nconst consumedCount = {};\nfor (const [oldText, newText] of replacements) {\n  const found = searchResults[oldText];\n  const n = consumedCount[oldText] || 0;\n  const range = found.items[n];\n  range.insertText(newText, \"Replace\");\n  consumedCount[oldText] = n + 1;\n}\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 two-digit / one-digit division\n# problems laid out in the 5-column table (problem rows 1, 5, 9, 13, 17;\n# 1-indexed), per the OOXML diff.\n$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\n  \"2026-01-06 Tuesday\", $false, $false, $false, $false, $false, $true,\n  1, $false, \"2026-01-07 Wednesday\", 2\n)\n\n$t = $d.Tables.Item(1)\n\n$values = @(\n  @(\"30\u00f75=\", \"98\u00f72=\", \"57\u00f78=\", \"31\u00f78=\", \"91\u00f72=\"),\n  @(\"42\u00f75=\", \"64\u00f77=\", \"59\u00f75=\", \"57\u00f74=\", \"59\u00f75=\"),\n  @(\"81\u00f72=\", \"92\u00f75=\", \"86\u00f77=\", \"78\u00f73=\", \"18\u00f74=\"),\n  @(\"88\u00f74=\", \"99\u00f74=\", \"17\u00f77=\", \"69\u00f79=\", \"35\u00f72=\"),\n  @(\"35\u00f74=\", \"54\u00f78=\", \"60\u00f78=\", \"98\u00f73=\", \"93\u00f76=\")\n)\n$rows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $rows.Length; $i++) {\n  $r = $rows[$i]\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($r, $c).Range.Text = $values[$i][$c - 1]\n  }\n}\n"}
